$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of data to insert right after the header row (before current row 2)
$newRows = @(
    @(-0.1327104717493057, 0.0224492978304624, 0.07544185966253281),
    @(-0.0704022198915481, 0.0503963828086853, -0.0980438739061355),
    @(0.0665843114256858, 0.0959058403968811, -0.0583376325666904),
    @(0.0452040284872055, 0.1351539343595504, -0.1539380401372909),
    @(0.2449569702148437, 0.4401284158229828, -0.2344195395708084)
)

# Insert 5 new rows starting at row 2, shifting existing data down
$insertRange = $ws.Range("A2:C6")
$insertRange.EntireRow.Insert()

# Populate the newly inserted rows with the new values
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $rowNum = 2 + $i
    $ws.Cells.Item($rowNum, 1).Value = $newRows[$i][0]
    $ws.Cells.Item($rowNum, 2).Value = $newRows[$i][1]
    $ws.Cells.Item($rowNum, 3).Value = $newRows[$i][2]
}

# Match the plain (unstyled) formatting of the other numeric data rows
$ws.Range("A2:C6").ClearFormats()

# The old rows 17-22 (now shifted down to rows 22-27 after the insert above)
# are no longer present in the final sheet; remove them so the sheet keeps
# only 20 data rows (A1:C21) as in the final dimension.
$ws.Range("A22:C27").EntireRow.Delete()
